# Refresh the Price (D) / Volume(1h) (E) columns and the handful of
# coin rows that swapped rank position (7<->8, 12<->13, 36<->37), as
# produced by the scheduled "Updated cryptos list" GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.506.82"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "'1.949.90"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'243.49"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'0.612"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'57.63"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "'0.374"
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("D10").Value = "'0.0786"
$ws.Range("E10").Value = "  -7.15%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "'2.236.60"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.826"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "'13.71"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "'20.98"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "'1.953.60"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "'36.424.13"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'69.37"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'0.0₃0845"
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").Value = "'228.25"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "'5.00"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").Value = "'9.08"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "'160.04"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "'0.134"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").Value = "'19.25"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").Value = "'0.0606"
$ws.Range("E33").Value = "  -5.06%  "
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'3.38"
$ws.Range("E36").Value = "  +11.82%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.23"
$ws.Range("E37").Value = "  +4.21%  "
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "'5.18"
$ws.Range("E39").Value = "  -15.73%  "
$ws.Range("D40").Value = "'0.0967"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'1.361.77"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "'15.67"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").Value = "'87.15"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'2.127.11"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'43.51"
$ws.Range("E51").Value = "  -1.98%  "
